{"js": "// Applies the program_writeup.docx edit described by the commit diff.\n// Strategy: locate each changed span of text with Body.search() (exact,\n// case-sensitive, whole phrase) and replace it in place with\n// Range.insertText(..., Word.InsertLocation.replace). Doing this as a\n// series of narrow, uniquely-matching replacements keeps formatting\n// (the run-level rPr with sz=24) intact because each replaced range sits\n// fully inside the original run(s) it overlaps.\n\nconst body = context.document.body;\n\nconst replacements = [\n  [\n    \"a single drive cpp file that handles\",\n    \"a single driver cpp file that handles\",\n  ],\n  [\n    \"rather handled the data loosely. \",\n    \"rather handled the data as it was retrieved. \",\n  ],\n  [\n    \"Once I switch to using the database I create the tables\",\n    \"Once I switch to using the Experiment database I create all of the tables\",\n  ],\n  [\n    \"application. From then I read in the user\\u2019s input and have a switch statement to determine what option the use\",\n    \"application. From there I prompt the user with a menu and read in the user\\u2019s input. I have a switch statement to determine what option the use\",\n  ],\n  [\n    \" selected. I broke out each of the tasks that the user can do into a separate function.\",\n    \" selected. Each of the tasks that the user can do has been broken into a separate function.\",\n  ],\n  [\n    \"They cannot break apart and insert just a single value into the database. \",\n    \"They cannot break the pieces apart and insert just a single result or parameter into the database. \",\n  ],\n  [\n    \"There isn\\u2019t any checking that needs to be done so it is straightforward to ask for and enter experiment meta data\",\n    \"There isn\\u2019t any checking that needs to be done when inserting an experiment, so I ask for and enter the experiment meta data\",\n  ],\n  [\n    \"run, it becomes more complicated. First,\",\n    \"run, it becomes more complicated since there are quite a few more checks that need to happen. First,\",\n  ],\n  [\n    \"Then I print out a list of all of the runs with an index and let the user select which run they\\u2019d like to display the information about. From there it is similar to printing out the meta data about the run, run parameters, and run results as it was for the experiment information.\",\n    \"Then I let the user select which run they\\u2019d like to display the information about. From there printing out the meta data about the run, run parameters, and run results is similar to printing the experiment meta data.\",\n  ],\n  [\n    \"parameters for the experiment they provide that may be aggregated, int or float.\",\n    \"parameters for the experiment they provided that may be aggregated, which are any that are an int or float.\",\n  ],\n  [\n    \"aggregated values which are calculated using SQL queries. Lastly for\",\n    \"aggregated values which are calculated using SQL queries. I make sure the user doesn\\u2019t enter the same date for the start and end date. Lastly for\",\n  ],\n  [\n    \"query to retrieve the results of each experiment with those values, by joining the Experiment and\",\n    \"query to retrieve the meta data of each experiment with that parameter, by joining the Experiment and\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + find);\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Applies the program_writeup.docx edit described by the commit diff.\n# Strategy: use Word's Find/Replace (Range.Find) to locate each changed\n# span of text (exact, case-sensitive, whole phrase) and replace it in\n# place. Each search string is unique in the document and sits fully\n# inside the original run(s) it overlaps, so formatting (the run-level\n# sz=24) carries over onto the replacement text automatically.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n\nReplace-Text \"a single drive cpp file that handles\" \"a single driver cpp file that handles\"\n\nReplace-Text \"rather handled the data loosely. \" \"rather handled the data as it was retrieved. \"\n\nReplace-Text \"Once I switch to using the database I create the tables\" \"Once I switch to using the Experiment database I create all of the tables\"\n\nReplace-Text \"application. From then I read in the user\u2019s input and have a switch statement to determine what option the use\" \"application. From there I prompt the user with a menu and read in the user\u2019s input. I have a switch statement to determine what option the use\"\n\nReplace-Text \" selected. I broke out each of the tasks that the user can do into a separate function.\" \" selected. Each of the tasks that the user can do has been broken into a separate function.\"\n\nReplace-Text \"They cannot break apart and insert just a single value into the database. \" \"They cannot break the pieces apart and insert just a single result or parameter into the database. \"\n\nReplace-Text \"There isn\u2019t any checking that needs to be done so it is straightforward to ask for and enter experiment meta data\" \"There isn\u2019t any checking that needs to be done when inserting an experiment, so I ask for and enter the experiment meta data\"\n\nReplace-Text \"run, it becomes more complicated. First,\" \"run, it becomes more complicated since there are quite a few more checks that need to happen. First,\"\n\nReplace-Text \"Then I print out a list of all of the runs with an index and let the user select which run they\u2019d like to display the information about. From there it is similar to printing out the meta data about the run, run parameters, and run results as it was for the experiment information.\" \"Then I let the user select which run they\u2019d like to display the information about. From there printing out the meta data about the run, run parameters, and run results is similar to printing the experiment meta data.\"\n\nReplace-Text \"parameters for the experiment they provide that may be aggregated, int or float.\" \"parameters for the experiment they provided that may be aggregated, which are any that are an int or float.\"\n\nReplace-Text \"aggregated values which are calculated using SQL queries. Lastly for\" \"aggregated values which are calculated using SQL queries. I make sure the user doesn\u2019t enter the same date for the start and end date. Lastly for\"\n\nReplace-Text \"query to retrieve the results of each experiment with those values, by joining the Experiment and\" \"query to retrieve the meta data of each experiment with that parameter, by joining the Experiment and\"\n"}
